$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 102, shifting existing rows 102:115 down to 103:116
$ws.Rows.Item(102).Insert()

# Fill in the new row 102 with the inserted data record
$ws.Cells.Item(102, 1).Value = 6
$ws.Cells.Item(102, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(102, 3).Value = "Metropolitana"
$ws.Cells.Item(102, 4).Value = 44984
$ws.Cells.Item(102, 5).Value = 13
$ws.Cells.Item(102, 6).Value = 100114007
$ws.Cells.Item(102, 7).Value = "Jengibre"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 285
$ws.Cells.Item(102, 11).Value = 17000
$ws.Cells.Item(102, 12).Value = 18000
$ws.Cells.Item(102, 13).Value = 17526
$ws.Cells.Item(102, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(102, 15).Value = "Perú"
$ws.Cells.Item(102, 16).Value = 1348
$ws.Cells.Item(102, 17).Value = 13
$ws.Cells.Item(102, 18).Value = "Hortaliza"
